# Auto-generated Excel COM-interop script
# Applies cached-value corrections to the FFXIV Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled-runner data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 41500
$ws.Range("J88").Value = 41500
$ws.Range("L88").Value = 41500
$ws.Range("N88").Value = -42312
$ws.Range("H91").Value = 41500
$ws.Range("J91").Value = 41500
$ws.Range("L91").Value = 41500
$ws.Range("N91").Value = -44308
$ws.Range("H101").Value = 2224.25
$ws.Range("J101").Value = 3494.75
$ws.Range("L101").Value = 10484.25
$ws.Range("N101").Value = -13728.25
$ws.Range("H132").Value = 50889.19
$ws.Range("I132").Value = 53239.65
$ws.Range("K132").Value = 159718.95
$ws.Range("M132").Value = -157188.95
$ws.Range("H137").Value = 1386.675
$ws.Range("I137").Value = 1320.8636
$ws.Range("J137").Value = 1467.1111
$ws.Range("K137").Value = 3962.5908
$ws.Range("L137").Value = 4401.3333
$ws.Range("M137").Value = -1412.5908
$ws.Range("N137").Value = -9501.3333
$ws.Range("H138").Value = 6546.4863
$ws.Range("I138").Value = 5323.1665
$ws.Range("J138").Value = 7133.68
$ws.Range("K138").Value = 15969.4995
$ws.Range("L138").Value = 21401.04
$ws.Range("M138").Value = -10829.4995
$ws.Range("N138").Value = -31681.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6039.357
$ws.Range("I2").Value = 6258.773
$ws.Range("J2").Value = 5234.8335
$ws.Range("K2").Value = 6258.773
$ws.Range("L2").Value = 5234.8335
$ws.Range("M2").Value = -6145.773
$ws.Range("N2").Value = -5460.8335
$ws.Range("H5").Value = 139.44444
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 5029.393
$ws.Range("I32").Value = 5030.4814
$ws.Range("K32").Value = 5030.4814
$ws.Range("M32").Value = -4743.4814
$ws.Range("H44").Value = 42000
$ws.Range("J44").Value = 42000
$ws.Range("L44").Value = 42000
$ws.Range("N44").Value = -42976
$ws.Range("H55").Value = 42000
$ws.Range("J55").Value = 42000
$ws.Range("L55").Value = 42000
$ws.Range("N55").Value = -42630
$ws.Range("H97").Value = 584
$ws.Range("I97").Value = 584
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 584
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -88
$ws.Range("N97").ClearContents()
$ws.Range("H116").Value = 6039.357
$ws.Range("I116").Value = 6258.773
$ws.Range("J116").Value = 5234.8335
$ws.Range("K116").Value = 6258.773
$ws.Range("L116").Value = 5234.8335
$ws.Range("M116").Value = -3964.773
$ws.Range("N116").Value = -9822.833500000001
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws.Range("H132").Value = 500000000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6039.357
$ws.Range("I3").Value = 6258.773
$ws.Range("J3").Value = 5234.8335
$ws.Range("K3").Value = 6258.773
$ws.Range("L3").Value = 5234.8335
$ws.Range("M3").Value = -6144.773
$ws.Range("N3").Value = -5462.8335
$ws.Range("H4").Value = 139.44444
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 1044.1
$ws.Range("I5").Value = 972.6875
$ws.Range("J5").Value = 1329.75
$ws.Range("K5").Value = 972.6875
$ws.Range("L5").Value = 1329.75
$ws.Range("M5").Value = -859.6875
$ws.Range("N5").Value = -1555.75
$ws.Range("H20").Value = 3045.2
$ws.Range("I20").Value = 2407.4167
$ws.Range("K20").Value = 2407.4167
$ws.Range("M20").Value = -2160.4167
$ws.Range("H22").Value = 376.2
$ws.Range("J22").Value = 376.6
$ws.Range("L22").Value = 376.6
$ws.Range("N22").Value = -722.6
$ws.Range("H99").Value = 4389.56
$ws.Range("I99").Value = 3955.5881
$ws.Range("K99").Value = 3955.5881
$ws.Range("M99").Value = -2457.5881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 669.87177
$ws.Range("J7").Value = 1073.1818
$ws.Range("L7").Value = 1073.1818
$ws.Range("N7").Value = -1299.1818
$ws.Range("H15").Value = 5999.4
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 5999.4
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5999.4
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -6339.4
$ws.Range("H22").Value = 835.25
$ws.Range("J22").Value = 1880.75
$ws.Range("L22").Value = 1880.75
$ws.Range("N22").Value = -2580.75
$ws.Range("H31").Value = 2562.3333
$ws.Range("J31").Value = 2671.818
$ws.Range("L31").Value = 2671.818
$ws.Range("N31").Value = -3261.818
$ws.Range("H34").Value = 2562.3333
$ws.Range("J34").Value = 2671.818
$ws.Range("L34").Value = 2671.818
$ws.Range("N34").Value = -3075.818
$ws.Range("H58").Value = 2491.75
$ws.Range("I58").Value = 2491.3845
$ws.Range("J58").Value = 2493.3333
$ws.Range("K58").Value = 2491.3845
$ws.Range("L58").Value = 2493.3333
$ws.Range("M58").Value = -2288.3845
$ws.Range("N58").Value = -2899.3333
$ws.Range("H132").Value = 2991.762
$ws.Range("I132").Value = 3284.5715
$ws.Range("K132").Value = 9853.7145
$ws.Range("M132").Value = -7323.7145
$ws.Range("H136").Value = 2491.75
$ws.Range("I136").Value = 2491.3845
$ws.Range("J136").Value = 2493.3333
$ws.Range("K136").Value = 7474.1535
$ws.Range("L136").Value = 7479.999899999999
$ws.Range("M136").Value = -4924.1535
$ws.Range("N136").Value = -12579.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 513.5714
$ws.Range("J12").Value = 510.8
$ws.Range("L12").Value = 1532.4
$ws.Range("N12").Value = -1878.4
$ws.Range("H103").Value = 420
$ws.Range("I103").Value = 420
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1260
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -381
$ws.Range("N103").ClearContents()
$ws.Range("H110").Value = 1998
$ws.Range("I110").Value = 1998
$ws.Range("K110").Value = 5994
$ws.Range("M110").Value = -1904
$ws.Range("H113").Value = 528.46155
$ws.Range("J113").Value = 548.5
$ws.Range("L113").Value = 1645.5
$ws.Range("N113").Value = -5985.5
$ws.Range("H121").Value = 52374.9
$ws.Range("I121").Value = 200223.8
$ws.Range("J121").Value = 3091.9333
$ws.Range("K121").Value = 600671.3999999999
$ws.Range("L121").Value = 9275.7999
$ws.Range("M121").Value = -599361.3999999999
$ws.Range("N121").Value = -11895.7999
$ws.Range("H124").Value = 1899
$ws.Range("I124").Value = 1899
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 5697
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -787
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 19037.4
$ws.Range("J4").Value = 25593.666
$ws.Range("L4").Value = 25593.666
$ws.Range("N4").Value = -25817.666
$ws.Range("H80").Value = 2988.6
$ws.Range("I80").Value = 2221
$ws.Range("K80").Value = 2221
$ws.Range("M80").Value = -1223
$ws.Range("H83").Value = 2988.6
$ws.Range("I83").Value = 2221
$ws.Range("K83").Value = 11105
$ws.Range("M83").Value = -6113
$ws.Range("H132").Value = 2998.6667
$ws.Range("I132").Value = 1998.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5995.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3465.5
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3636.4614
$ws.Range("I22").Value = 1328.4286
$ws.Range("J22").Value = 6329.1665
$ws.Range("K22").Value = 1328.4286
$ws.Range("L22").Value = 6329.1665
$ws.Range("M22").Value = -1033.4286
$ws.Range("N22").Value = -6919.1665
$ws.Range("H27").Value = 3636.4614
$ws.Range("I27").Value = 1328.4286
$ws.Range("J27").Value = 6329.1665
$ws.Range("K27").Value = 1328.4286
$ws.Range("L27").Value = 6329.1665
$ws.Range("M27").Value = -1221.4286
$ws.Range("N27").Value = -6543.1665
$ws.Range("H100").Value = 1172.6
$ws.Range("I100").Value = 1149
$ws.Range("K100").Value = 1149
$ws.Range("M100").Value = -608
$ws.Range("H136").Value = 5380.1333
$ws.Range("I136").Value = 5101.4
$ws.Range("K136").Value = 15304.2
$ws.Range("M136").Value = -12754.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 952.84
$ws.Range("I113").Value = 578.3570999999999
$ws.Range("K113").Value = 1735.0713
$ws.Range("M113").Value = 434.9287000000002

Write-Host "Applied cached-value updates across 8 sheets."
